$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row on this sheet
$lastRow = $ws.UsedRange.Rows.Count

# 1) Update column C (Förändrad) from 45184 to 45186 for every data row (2..lastRow)
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45186
    }
}

# 2) Add a second argument (the "Beteckning" text in column A) to every
#    HYPERLINK formula that doesn't already have one, in columns S..Y
$hlCols = @(19, 20, 21, 22, 23, 24, 25)  # S, T, U, V, W, X, Y

for ($r = 2; $r -le $lastRow; $r++) {
    $beteckning = $ws.Cells.Item($r, 1).Value2
    if ([string]::IsNullOrEmpty($beteckning)) { continue }

    foreach ($c in $hlCols) {
        $cell = $ws.Cells.Item($r, $c)
        $f = $cell.Formula
        if ($f -and $f -like '*HYPERLINK(*' -and $f -notlike '*,*') {
            $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $beteckning + '")'
            $cell.Formula = $newFormula
        }
    }
}
